$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column B (RF descriptions, rows 5-15) lost the extra "applyFont"
#    style (cellXfs index 4) and now share the plain bordered style
#    already used by column A. Re-apply that format (copy formats only)
#    so the style collapses onto the same xf as A5:A15 / B2:B4.
# ---------------------------------------------------------------------
foreach ($r in 5..15) {
    $ws.Range("A$r").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. New non-functional requirements RNF07-RNF14 fill F8:I15, mirroring
#    the existing F2:I7 block (format copied from the row above it).
# ---------------------------------------------------------------------
foreach ($r in 8..15) {
    $ws.Range("F7:I7").Copy() | Out-Null
    $ws.Range("F$r" + ":I$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

$ws.Range("F8").Value = "RNF07"
$ws.Range("F9").Value = "RNF08"
$ws.Range("F10").Value = "RNF09"
$ws.Range("F11").Value = "RNF10"
$ws.Range("F12").Value = "RNF11"
$ws.Range("F13").Value = "RNF12"
$ws.Range("F14").Value = "RNF13"
$ws.Range("F15").Value = "RNF14"

$ws.Range("G8").Value = "banco de dados de temperatura "
$ws.Range("G9").Value = "banco de dados usário"
$ws.Range("G10").Value = "sistema fazer backup periodicamente"
$ws.Range("G11").Value = "guardar backup em um local seguro"
$ws.Range("G12").Value = "verificar se a bugs periodiacamente"
$ws.Range("G13").Value = "verificar se a feedback ou sugestões dos usuários"
# G14 / G15 stay blank (only formatted, like the template row)

$ws.Range("H8").Value = "importante"
$ws.Range("H9").Value = "importante"
$ws.Range("H10").Value = "essencial"
$ws.Range("H11").Value = "essencial"
$ws.Range("H12").Value = "importante"
$ws.Range("H13").Value = "desejavel"
$ws.Range("H14").Value = "essencial"
$ws.Range("H15").Value = "essencial"

# I8:I15 remain empty (date-formatted placeholders), matching I2:I7's look.

# ---------------------------------------------------------------------
# 3. Cosmetic follow-ups from the resave: widen column G a bit and move
#    the live selection to H13.
# ---------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 47.25

$ws.Range("H13").Select() | Out-Null

Write-Output "requisitos nao funcionais atualizados"
